$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "77÷5=15, 2"
$t.Cell(1,2).Range.Text = "84÷2=42, 0"
$t.Cell(1,3).Range.Text = "67÷5=13, 2"
$t.Cell(1,4).Range.Text = "74÷5=14, 4"
$t.Cell(1,5).Range.Text = "94÷2=47, 0"
$t.Cell(5,1).Range.Text = "94÷8=11, 6"
$t.Cell(5,2).Range.Text = "14÷9=1, 5"
$t.Cell(5,3).Range.Text = "19÷4=4, 3"
$t.Cell(5,4).Range.Text = "17÷9=1, 8"
$t.Cell(5,5).Range.Text = "23÷7=3, 2"
$t.Cell(9,1).Range.Text = "49÷9=5, 4"
$t.Cell(9,2).Range.Text = "10÷7=1, 3"
$t.Cell(9,3).Range.Text = "42÷2=21, 0"
$t.Cell(9,4).Range.Text = "17÷2=8, 1"
$t.Cell(9,5).Range.Text = "92÷6=15, 2"
$t.Cell(13,1).Range.Text = "93÷5=18, 3"
$t.Cell(13,2).Range.Text = "92÷3=30, 2"
$t.Cell(13,3).Range.Text = "42÷2=21, 0"
$t.Cell(13,4).Range.Text = "17÷2=8, 1"
$t.Cell(13,5).Range.Text = "59÷6=9, 5"
$t.Cell(17,1).Range.Text = "14÷6=2, 2"
$t.Cell(17,2).Range.Text = "42÷2=21, 0"
$t.Cell(17,3).Range.Text = "53÷5=10, 3"
$t.Cell(17,4).Range.Text = "26÷8=3, 2"
$t.Cell(17,5).Range.Text = "22÷9=2, 4"
